$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='56.935.09'; E='  +4.49%  '}
    @{Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='3.239.61'; E='  +1.83%  '}
    @{Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='1.00'; E='  +0.16%  '}
    @{Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='395.10'; E='  -2.18%  '}
    @{Row=6; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='108.01'; E='  -1.09%  '}
    @{Row=7; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.582'; E='  +5.71%  '}
    @{Row=8; B='LidoStakedEther'; C='https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'; D='3.232.93'; E='  +1.84%  '}
    @{Row=9; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='1.00'; E='  +0.04%  '}
    @{Row=10; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.620'; E='  -0.10%  '}
    @{Row=11; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='38.97'; E='  +0.00%  '}
    @{Row=12; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.0994'; E='  +12.57%  '}
    @{Row=13; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.142'; E='  +1.69%  '}
    @{Row=14; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='3.757.99'; E='  +2.09%  '}
    @{Row=15; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='8.23'; E='  +2.07%  '}
    @{Row=16; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='19.04'; E='  -0.97%  '}
    @{Row=17; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='3.241.19'; E='  +1.69%  '}
    @{Row=18; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='1.03'; E='  -3.32%  '}
    @{Row=19; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='10.70'; E='  +2.05%  '}
    @{Row=20; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='56.849.19'; E='  +4.64%  '}
    @{Row=21; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='3.32'; E='  +0.54%  '}
    @{Row=22; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.0000110'; E='  +10.87%  '}
    @{Row=23; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='12.83'; E='  -0.36%  '}
    @{Row=24; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='294.75'; E='  +7.42%  '}
    @{Row=25; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='74.08'; E='  +3.43%  '}
    @{Row=26; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='3.16'; E='  -4.00%  '}
    @{Row=27; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='28.00'; E='  +0.88%  '}
    @{Row=28; B='LEO'; C='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; D='4.35'; E='  +1.27%  '}
    @{Row=29; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='7.56'; E='  -5.66%  '}
    @{Row=30; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='7.19'; E='  -2.64%  '}
    @{Row=31; B='Kaspa'; C='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D='0.167'; E='  -1.54%  '}
    @{Row=32; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.00'; E='  -0.04%  '}
    @{Row=33; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='11.24'; E='  +1.39%  '}
    @{Row=34; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.108'; E='  -3.85%  '}
    @{Row=35; B='InjectiveProtocol'; C='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D='38.94'; E='  +5.67%  '}
    @{Row=36; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.0479'; E='  -3.95%  '}
    @{Row=37; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='2.13'; E='  +2.03%  '}
    @{Row=38; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='51.35'; E='  +1.28%  '}
    @{Row=39; B='FirstDigitalUSD'; C='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D='1.00'; E='  +0.23%  '}
    @{Row=40; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='3.46'; E='  -5.63%  '}
    @{Row=41; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='2.88'; E='  +1.64%  '}
    @{Row=42; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='135.13'; E='  +3.76%  '}
    @{Row=43; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.122'; E='  +4.11%  '}
    @{Row=44; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.88'; E='  -2.80%  '}
    @{Row=45; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='3.93'; E='  -4.26%  '}
    @{Row=46; B='Celestia'; C='https://coinranking.com/coin/YQcD0lBl7+celestia-tia'; D='16.91'; E='  -2.25%  '}
    @{Row=47; B='TheGraph'; C='https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; D='0.278'; E='  -4.56%  '}
    @{Row=48; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='22.02'; E='  -0.84%  '}
    @{Row=49; B='WEMIXToken'; C='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D='2.13'; E='  +3.08%  '}
    @{Row=50; B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='2.144.80'; E='  +2.76%  '}
    @{Row=51; B='ApeXProtocol'; C='https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'; D='2.33'; E='  -6.92%  '}
)

foreach ($item in $data) {
    $r = $item.Row

    $cB = $ws.Cells.Item($r, 2)
    $cB.NumberFormat = "@"
    $cB.Value2 = $item.B
    $cB.Style = "Normal"

    $cC = $ws.Cells.Item($r, 3)
    $cC.NumberFormat = "@"
    $cC.Value2 = $item.C
    $cC.Style = "Normal"

    $cD = $ws.Cells.Item($r, 4)
    $cD.NumberFormat = "@"
    $cD.Value2 = $item.D
    $cD.Style = "Normal"

    $cE = $ws.Cells.Item($r, 5)
    $cE.NumberFormat = "@"
    $cE.Value2 = $item.E
    $cE.Style = "Normal"
}

$wb.Save()